$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: Panel B (E-mini Futures) - Ann Window Volume
$ws.Range("D26").Value = 4178.681720430108
$ws.Range("E26").Value = 2919.149727889589
$ws.Range("F26").Value = 1727.225806451613
$ws.Range("G26").Value = 3726.612903225806
$ws.Range("H26").Value = 5975.967741935484
$ws.Range("I26").Value = 45
$ws.Range("J26").Value = 4467.179234972678
$ws.Range("K26").Value = 1841.843065169242
$ws.Range("L26").Value = 3168.852459016393
$ws.Range("M26").Value = 3878.131147540983
$ws.Range("N26").Value = 5485.163934426229
$ws.Range("O26").Value = 45
$ws.Range("P26").Value = 4188.48668503214
$ws.Range("Q26").Value = 1506.06581302919
$ws.Range("R26").Value = 2909.413223140496
$ws.Range("S26").Value = 4010.363636363636
$ws.Range("T26").Value = 5059.297520661157
$ws.Range("U26").Value = 45
$ws.Range("V26").Value = 3771.793121693122
$ws.Range("W26").Value = 1222.595510054101
$ws.Range("X26").Value = 2936.738095238095
$ws.Range("Y26").Value = 3673.92380952381
$ws.Range("Z26").Value = 4499.547619047619
$ws.Range("AA26").Value = 45
$ws.Range("AB26").Value = 1101.638114478115
$ws.Range("AC26").Value = 262.6569382038755
$ws.Range("AD26").Value = 940.1931818181819
$ws.Range("AE26").Value = 1053.273484848485
$ws.Range("AF26").Value = 1267.327272727273
$ws.Range("AG26").Value = 45

# Row 27: Panel B (E-mini Futures) - Diff (Ann - Non)
$ws.Range("D27").Value = 1869.684050179211
$ws.Range("J27").Value = 2080.352367941713
$ws.Range("P27").Value = 1797.004545454545
$ws.Range("V27").Value = 1234.992328042328
$ws.Range("AB27").Value = 107.3329419191919

# Row 28: Panel B (E-mini Futures) - # Obs
$ws.Range("D28").Value = 45
$ws.Range("J28").Value = 45
$ws.Range("P28").Value = 45
$ws.Range("V28").Value = 45
$ws.Range("AB28").Value = 45
